# Refresh the Price (column D) and Volume(1h) (column E) figures in the
# cryptos list, as produced by the scheduled GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.441.64'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.885.09'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  -0.74%  '
$__s = $ws.Range("D5").Style
$ws.Range("D5").Formula = "'246.65"
$ws.Range("D5").Style = $__s
$ws.Range("E5").Value = '  -3.78%  '
$__s = $ws.Range("D6").Style
$ws.Range("D6").Formula = "'0.688"
$ws.Range("D6").Style = $__s
$ws.Range("E6").Value = '  -5.80%  '
$__s = $ws.Range("D7").Style
$ws.Range("D7").Formula = "'0.999"
$ws.Range("D7").Style = $__s
$ws.Range("E7").Value = '  -0.82%  '
$__s = $ws.Range("D8").Style
$ws.Range("D8").Formula = "'43.35"
$ws.Range("D8").Style = $__s
$ws.Range("E8").Value = '  +4.60%  '
$__s = $ws.Range("D9").Style
$ws.Range("D9").Formula = "'0.352"
$ws.Range("D9").Style = $__s
$ws.Range("E9").Value = '  -4.41%  '
$__s = $ws.Range("D10").Style
$ws.Range("D10").Formula = "'0.0737"
$ws.Range("D10").Style = $__s
$ws.Range("E10").Value = '  -3.32%  '
$__s = $ws.Range("D11").Style
$ws.Range("D11").Formula = "'0.0969"
$ws.Range("D11").Style = $__s
$ws.Range("E11").Value = '  -2.00%  '
$__s = $ws.Range("D12").Style
$ws.Range("D12").Formula = "'13.04"
$ws.Range("D12").Style = $__s
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '2.152.28'
$ws.Range("E13").Value = '  -1.73%  '
$__s = $ws.Range("D14").Style
$ws.Range("D14").Formula = "'0.739"
$ws.Range("D14").Style = $__s
$ws.Range("E14").Value = '  +0.40%  '
$__s = $ws.Range("D15").Style
$ws.Range("D15").Formula = "'4.94"
$ws.Range("D15").Style = $__s
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").Value = '1.865.01'
$ws.Range("E16").Value = '  -3.58%  '
$ws.Range("D17").Value = '35.350.88'
$ws.Range("E17").Value = '  -0.03%  '
$__s = $ws.Range("D18").Style
$ws.Range("D18").Formula = "'73.54"
$ws.Range("D18").Style = $__s
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").Value = '0.0₃0821'
$ws.Range("E19").Value = '  -3.45%  '
$__s = $ws.Range("D20").Style
$ws.Range("D20").Formula = "'246.25"
$ws.Range("D20").Style = $__s
$ws.Range("E20").Value = '  +0.31%  '
$__s = $ws.Range("D21").Style
$ws.Range("D21").Formula = "'12.81"
$ws.Range("D21").Style = $__s
$ws.Range("E21").Value = '  -2.53%  '
$__s = $ws.Range("D22").Style
$ws.Range("D22").Formula = "'4.93"
$ws.Range("D22").Style = $__s
$ws.Range("E22").Value = '  -4.60%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("E25").Value = '  -11.37%  '
$__s = $ws.Range("D26").Style
$ws.Range("D26").Formula = "'165.55"
$ws.Range("D26").Style = $__s
$ws.Range("E26").Value = '  -0.66%  '
$__s = $ws.Range("D27").Style
$ws.Range("D27").Formula = "'8.47"
$ws.Range("D27").Style = $__s
$ws.Range("E27").Value = '  -3.69%  '
$__s = $ws.Range("D28").Style
$ws.Range("D28").Formula = "'18.29"
$ws.Range("D28").Style = $__s
$ws.Range("E28").Value = '  -2.96%  '
$ws.Range("E29").Value = '  -4.76%  '
$ws.Range("D30").Value = '4.128.42'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  +4.75%  '
$__s = $ws.Range("D32").Style
$ws.Range("D32").Formula = "'4.24"
$ws.Range("D32").Style = $__s
$ws.Range("E32").Value = '  -3.40%  '
$__s = $ws.Range("D33").Style
$ws.Range("D33").Formula = "'0.0578"
$ws.Range("D33").Style = $__s
$ws.Range("E33").Value = '  -2.48%  '
$__s = $ws.Range("D34").Style
$ws.Range("D34").Formula = "'4.18"
$ws.Range("D34").Style = $__s
$ws.Range("E34").Value = '  -2.44%  '
$ws.Range("E35").Value = '  -0.80%  '
$__s = $ws.Range("D36").Style
$ws.Range("D36").Formula = "'0.849"
$ws.Range("D36").Style = $__s
$ws.Range("E36").Value = '  -7.62%  '
$__s = $ws.Range("D37").Style
$ws.Range("D37").Formula = "'1.98"
$ws.Range("D37").Style = $__s
$ws.Range("E37").Value = '  -3.36%  '
$ws.Range("E38").Value = '  -21.32%  '
$__s = $ws.Range("D39").Style
$ws.Range("D39").Formula = "'0.0692"
$ws.Range("D39").Style = $__s
$ws.Range("E39").Value = '  +7.41%  '
$__s = $ws.Range("D40").Style
$ws.Range("D40").Formula = "'97.02"
$ws.Range("D40").Style = $__s
$ws.Range("E40").Value = '  -0.97%  '
$__s = $ws.Range("D41").Style
$ws.Range("D41").Formula = "'16.95"
$ws.Range("D41").Style = $__s
$ws.Range("E41").Value = '  -1.21%  '
$__s = $ws.Range("D42").Style
$ws.Range("D42").Formula = "'0.0215"
$ws.Range("D42").Style = $__s
$ws.Range("E42").Value = '  -2.81%  '
$ws.Range("E43").Value = '  -4.41%  '
$ws.Range("D44").Value = '1.287.50'
$ws.Range("E44").Value = '  -4.40%  '
$__s = $ws.Range("D45").Style
$ws.Range("D45").Formula = "'2.34"
$ws.Range("D45").Style = $__s
$ws.Range("E45").Value = '  -5.78%  '
$__s = $ws.Range("D46").Style
$ws.Range("D46").Formula = "'0.0811"
$ws.Range("D46").Style = $__s
$ws.Range("E46").Value = '  +7.34%  '
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("E48").Value = '  -0.93%  '
$__s = $ws.Range("D49").Style
$ws.Range("D49").Formula = "'12.10"
$ws.Range("D49").Style = $__s
$ws.Range("E49").Value = '  +3.00%  '
$__s = $ws.Range("D50").Style
$ws.Range("D50").Formula = "'43.16"
$ws.Range("D50").Style = $__s
$ws.Range("E50").Value = '  -4.73%  '
$__s = $ws.Range("D51").Style
$ws.Range("D51").Formula = "'6.24"
$ws.Range("D51").Style = $__s
$ws.Range("E51").Value = '  -7.53%  '
